# UCBI_YR_FIN.xlsx - "Doing Updates for Financials"
# A new (most recent) reporting period is inserted as column D on sheet "UCBI",
# pushing the previously existing columns D:K one column to the right (E:L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column D; existing D:K data (and formatting)
# shifts to E:L automatically.
$ws.Columns("D:D").Insert()

# The freshly inserted column D has default formatting - clone the number
# formatting/style from column E (which now holds what used to be in D) onto
# the new column D for every row that actually contains data, so the new
# column looks consistent with the rest of the table (date format row 7/38/80,
# plain-number style elsewhere).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest fiscal-period figures.

# --- Income Statement ---
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 500100
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = -6800
$ws.Range("D17").Value = 70800
$ws.Range("D18").Value = 429300
$ws.Range("D20").Value = -213300
$ws.Range("D21").Value = 246900
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 215900
$ws.Range("D24").Value = 49800
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 166100
$ws.Range("D27").Value = 164900
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 213300
$ws.Range("D33").Value = 164900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 164900

# --- Balance Sheet ---
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 126100
$ws.Range("D42").Value = 225900
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 206100
$ws.Range("D49").Value = 324100
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 64200
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 12573200
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 127500
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 267200
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 11115600
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -79700
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1457600
$ws.Range("D77").Value = 0

# --- Cash Flow Statement ---
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 164900
$ws.Range("D83").Value = 31000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 270000
$ws.Range("D91").Value = -17600
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -359500
$ws.Range("D96").Value = -41600
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 102500
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 13000

$wb.Save()
